# Applies the "fixes and word updates" commit:
#  - Sheet "4.2": a batch of C-column cells change from "wrong" to "correct"
#  - Sheet "4.3": a batch of C-column cells change from "wrong" to "correct"
#  - Sheet "4.5": a new C column ("correct") is added for every existing row
#  - Sheet "4.6": a new C column ("correct") is added for every existing row,
#                 plus two French typo fixes in column A
#  - The workbook's active tab changes to sheet "4.6" (index 30, 0-based)

$wb = $excel.ActiveWorkbook

# --- Sheet "4.2": rows whose column C currently reads "wrong" -> "correct"
$ws42 = $wb.Worksheets.Item("4.2")
$rows42 = @(4,5,6,7,8,11,15,16,18,20,23,25,31,36,39,41,42,48,49,51)
foreach ($r in $rows42) {
    $ws42.Cells.Item($r, 3).Value = "correct"
}

# --- Sheet "4.3": rows whose column C currently reads "wrong" -> "correct"
$ws43 = $wb.Worksheets.Item("4.3")
$rows43 = @(3,4,6,10,11,16,17,18,19,21,23,27,28,29,30,32,33,36,38,39,40,42,45,46,47,48,51,52)
foreach ($r in $rows43) {
    $ws43.Cells.Item($r, 3).Value = "correct"
}

# --- Sheet "4.5": add column C = "correct" for all 57 existing rows
$ws45 = $wb.Worksheets.Item("4.5")
for ($r = 1; $r -le 57; $r++) {
    $ws45.Cells.Item($r, 3).Value = "correct"
}

# --- Sheet "4.6": add column C = "correct" for all 60 existing rows
$ws46 = $wb.Worksheets.Item("4.6")
for ($r = 1; $r -le 60; $r++) {
    $ws46.Cells.Item($r, 3).Value = "correct"
}

# Fix two French typos in sheet "4.6" column A
$ws46.Range("A44").Value = "se rassembler"
$ws46.Range("A47").Value = "un repas"

# Make "4.6" the active tab, matching the saved workbook view state
$ws46.Activate()
